$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Veiculos")

$ws.Range("A5").Value = -1
$ws.Range("B5").Value = "QWE1R12"
$ws.Range("C5").Value = "DASDA"
$ws.Range("D5").Value = "Chevrolet"
$ws.Range("E5").Value = 2000
$ws.Range("F5").Value = "Ativo"

$ws.Range("A6").Value = 55
$ws.Range("B6").Value = "QWE2E12"
$ws.Range("C6").Value = "ADA"
$ws.Range("D6").Value = "Chevrolet"
$ws.Range("E6").Value = 2025
$ws.Range("F6").Value = "Ativo"
